$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6 (GraphSMOTE 21') ---
$ws.Range("B6").Value = "80.030 (0.646)"
$ws.Range("C6").Value = "76.116 (0.889)"
$ws.Range("D6").Value = "76.418 (1.255)"
$ws.Range("E6").Value = "84.970 (0.527)"
$ws.Range("F6").Value = "82.680 (0.728)"
$ws.Range("G6").Value = "83.301 (0.667)"
$ws.Range("H6").Value = "64.710 (0.803)"
$ws.Range("I6").Value = "58.119 (0.733)"
$ws.Range("J6").Value = "51.492 (1.577)"
$ws.Range("K6").Value = "70.960 (0.450)"
$ws.Range("L6").Value = "64.043 (0.501)"
$ws.Range("M6").Value = "60.787 (0.863)"
$ws.Range("N6:S6").ClearContents()

# --- Row 12 (GraphSHA 23') add new cells N12:S12 ---
$ws.Range("N12").Value = "67.900 (2.106)"
$ws.Range("O12").Value = "65.440 (3.259)"
$ws.Range("P12").Value = "64.121 (2.769)"
$ws.Range("Q12").Value = "80.200 (0.346)"
$ws.Range("R12").Value = "77.456 (0.460)"
$ws.Range("S12").Value = "77.953 (0.362)"

# --- Row 21 (Vanilla (GCN)) ---
$ws.Range("B21").Value = "90.933 (0.159)"
$ws.Range("C21").Value = "83.850 (0.553)"
$ws.Range("D21").Value = "83.799 (0.841)"
$ws.Range("H21").Value = "72.892 (3.246)"
$ws.Range("I21").Value = "36.448 (7.049)"
$ws.Range("J21").Value = "33.942 (7.430)"
$ws.Range("N21").Value = "92.796 (0.091)"
$ws.Range("O21").Value = "88.617 (0.233)"
$ws.Range("P21").Value = "89.867 (0.179)"

# --- Row 23 (GraphSMOTE 21') new cells ---
$ws.Range("B23").Value = "91.534 (0.274)"
$ws.Range("C23").Value = "86.850 (0.520)"
$ws.Range("D23").Value = "88.133 (0.511)"
$ws.Range("H23").Value = "71.735 (1.678)"
$ws.Range("I23").Value = "40.059 (1.465)"
$ws.Range("J23").Value = "37.582 (1.500)"
$ws.Range("N23").Value = "92.450 (0.084)"
$ws.Range("O23").Value = "87.804 (0.243)"
$ws.Range("P23").Value = "89.186 (0.194)"

# --- Row 25 (GraphENS 22') ---
$ws.Range("B25").Value = "90.961 (0.396)"
$ws.Range("C25").Value = "91.290 (0.150)"
$ws.Range("D25").Value = "88.617 (0.373)"
$ws.Range("H25").Value = "82.284 (0.812)"
$ws.Range("I25").Value = "87.639 (0.231)"
$ws.Range("J25").Value = "79.302 (0.672)"
$ws.Range("N25").Value = "92.172 (0.225)"
$ws.Range("O25").Value = "90.492 (0.177)"
$ws.Range("P25").Value = "85.814 (3.212)"

# --- Row 29 (GraphSHA 23') ---
$ws.Range("B29").Value = "92.112 (0.176)"
$ws.Range("C29").Value = "91.075 (0.149)"
$ws.Range("D29").Value = "89.862 (0.165)"
$ws.Range("H29").Value = "84.051 (0.403)"
$ws.Range("I29").Value = "84.767 (0.537)"
$ws.Range("J29").Value = "81.291 (0.436)"
$ws.Range("N29").Value = "92.533 (0.060)"
$ws.Range("O29").Value = "90.140 (0.143)"
$ws.Range("P29").Value = "84.724 (2.489)"
